$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the prep material for week 8 to D9 (Vorbereitung column)
$ws.Range("D9").Value = "prep/p08.html"

# Add the slides link for session 06 to E7 (Folien column, week 6 row)
$ws.Range("E7").Value = "slides/slides.html#/sitzung-06-warum-diese-medienerlebnisse-manchmal-trotzdem-unterhaltsam-sein-können"

# Update the selected cell to reflect the last edited cell
$ws.Range("E7").Select()
